$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Ukuran" header in F1 to "Tipe"
$ws.Range("F1").Value = "Tipe"

# F1 should now use the same header style as the other plain headers
# (centered, bold, filled) instead of its previous unaligned variant.
$ws.Range("F1").HorizontalAlignment = -4108 ; # xlCenter

# Column F becomes a bit wider and its body cells get center alignment
$ws.Range("F1").EntireColumn.ColumnWidth = 10.6640625
$ws.Range("F1").EntireColumn.HorizontalAlignment = -4108 ; # xlCenter

# Move the active selection to F3
$ws.Range("F3").Select()
